$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = '[["CDU updated", 0.03370302835384297], ["Steel cylinder", 2.4821260976820004e-05], ["nitrous oxide production", 0.010968005818038482], ["Use phase CDU", 0.001199576119063503], ["EoL CDU", -0.0012243440526300189], ["EoL 50L cylinder - FU", 0.0015934554488819326], ["avoided heat", -2.6182230591314203e-09], ["avoided electricity", -2.5724139875873146e-08]]'
$ws.Range("A4").Value = '[["MDU updated", 0.06909101812703686], ["Steel cylinder", 2.4821260976820004e-05], ["nitrous oxide production", 0.010968005818038482], ["Use phase MDU", 0.00045196796336067954], ["EoL MDU", -0.007055596602464493], ["EoL 50L cylinder - FU", 0.0015934554488819326], ["avoided heat", -7.706244714779625e-09], ["avoided electricity", -1.2379591978014953e-07]]'
$ws.Range("B3").Value = '[["CDU updated", 0.12862424689466137], ["Steel cylinder", 0.005456229514756893], ["nitrous oxide production", 2.7735767368862367], ["Use phase CDU", 43.44535637749805], ["EoL CDU", 0.00028463836607178506], ["EoL 50L cylinder - FU", 1.37975498672541], ["avoided heat", -1.1022661664048948e-06], ["avoided electricity", -1.003541297311142e-05]]'
$ws.Range("B4").Value = '[["MDU updated", 0.3591249304896528], ["Steel cylinder", 0.005456229514756893], ["nitrous oxide production", 2.7735767368862367], ["Use phase MDU", 43.40432983037498], ["EoL MDU", 0.00705331734479795], ["EoL 50L cylinder - FU", 1.37975498672541], ["avoided heat", -3.2443121259512005e-06], ["avoided electricity", -4.8294838442597295e-05]]'
$ws.Range("C3").Value = '[["CDU updated", 1.4827584819565856], ["Steel cylinder", 0.02689681145503979], ["nitrous oxide production", 15.635308294302359], ["Use phase CDU", 0.7412027383344805], ["EoL CDU", 0.004577101780804119], ["EoL 50L cylinder - FU", 2.596880646960822], ["avoided heat", -1.430092102445719e-06], ["avoided electricity", -2.0189326346486028e-05]]'
$ws.Range("C4").Value = '[["MDU updated", 5.562379935789622], ["Steel cylinder", 0.02689681145503979], ["nitrous oxide production", 15.635308294302359], ["Use phase MDU", 0.27926522273878285], ["EoL MDU", 0.10369155230631812], ["EoL 50L cylinder - FU", 2.596880646960822], ["avoided heat", -4.209205807603736e-06], ["avoided electricity", -9.71599531360455e-05]]'
$ws.Range("D3").Value = '[["CDU updated", 2.0996855665639584], ["Steel cylinder", 0.06471586600422786], ["nitrous oxide production", 46.31028189903245], ["Use phase CDU", 0.6692525252014253], ["EoL CDU", -0.0026206906884051256], ["EoL 50L cylinder - FU", 6.953517004760155], ["avoided heat", -1.3743984342168013e-05], ["avoided electricity", -0.0001215476478635583]]'
$ws.Range("D4").Value = '[["MDU updated", 5.548337101286209], ["Steel cylinder", 0.06471586600422786], ["nitrous oxide production", 46.31028189903245], ["Use phase MDU", 0.25215632087225787], ["EoL MDU", 0.03362627844392727], ["EoL 50L cylinder - FU", 6.953517004760155], ["avoided heat", -4.045282021638844e-05], ["avoided electricity", -0.0005849409518448048]]'
$ws.Range("E3").Value = '[["CDU updated", 0.00010984081423254179], ["Steel cylinder", 4.164145755153704e-06], ["nitrous oxide production", 0.0009169092173445685], ["Use phase CDU", 2.849641883709527e-05], ["EoL CDU", 6.532498571501736e-07], ["EoL 50L cylinder - FU", 0.00012311536065306782], ["avoided heat", -1.663567245821682e-10], ["avoided electricity", -3.937282608610133e-09]]'
$ws.Range("E4").Value = '[["MDU updated", 0.00036637421639992547], ["Steel cylinder", 4.164145755153704e-06], ["nitrous oxide production", 0.0009169092173445685], ["Use phase MDU", 1.073668288339174e-05], ["EoL MDU", 9.471947961735349e-06], ["EoL 50L cylinder - FU", 0.00012311536065306782], ["avoided heat", -4.896395763948186e-10], ["avoided electricity", -1.8947942450911875e-08]]'
$ws.Range("F3").Value = '[["CDU updated", 0.00012550417287013446], ["Steel cylinder", 5.439366377934605e-06], ["nitrous oxide production", 0.005230909219588117], ["Use phase CDU", 0.00027925130314251294], ["EoL CDU", 6.228837949449939e-06], ["EoL 50L cylinder - FU", 0.00065373948486032], ["avoided heat", -7.992765374258922e-10], ["avoided electricity", -6.55656546450892e-09]]'
$ws.Range("F4").Value = '[["MDU updated", 0.00036185641280330627], ["Steel cylinder", 5.439366377934605e-06], ["nitrous oxide production", 0.005230909219588117], ["Use phase MDU", 0.00010521436759316886], ["EoL MDU", 4.849556449600232e-05], ["EoL 50L cylinder - FU", 0.00065373948486032], ["avoided heat", -2.352519419891193e-09], ["avoided electricity", -3.155308811856111e-08]]'
$ws.Range("G3").Value = '[["CDU updated", 0.001298493755320083], ["Steel cylinder", 5.1440088432858044e-05], ["nitrous oxide production", 0.029899230463375122], ["Use phase CDU", 0.004351892777417035], ["EoL CDU", 9.288837439875151e-05], ["EoL 50L cylinder - FU", 0.004205256079904196], ["avoided heat", -8.797126803003337e-09], ["avoided electricity", -6.97471342106168e-08]]'
$ws.Range("G4").Value = '[["MDU updated", 0.0037610866719553176], ["Steel cylinder", 5.1440088432858044e-05], ["nitrous oxide production", 0.029899230463375122], ["Use phase MDU", 0.0016396759522928257], ["EoL MDU", 0.0007148061559304394], ["EoL 50L cylinder - FU", 0.004205256079904196], ["avoided heat", -2.589268003529705e-08], ["avoided electricity", -3.3565400722030797e-07]]'
$ws.Range("H3").Value = '[["CDU updated", 3.2375682748638255e-10], ["Steel cylinder", 8.120702205798173e-11], ["nitrous oxide production", 1.3815806307562065e-09], ["Use phase CDU", 3.112433080852018e-10], ["EoL CDU", 7.628871961180972e-12], ["EoL 50L cylinder - FU", 2.454274504847323e-10], ["avoided heat", -6.797015971576207e-16], ["avoided electricity", -4.919358458857522e-15]]'
$ws.Range("H4").Value = '[["MDU updated", 1.0716849152455745e-09], ["Steel cylinder", 8.120702205798173e-11], ["nitrous oxide production", 1.3815806307562065e-09], ["Use phase MDU", 1.1726809314504075e-10], ["EoL MDU", 8.149709496030253e-11], ["EoL 50L cylinder - FU", 2.454274504847323e-10], ["avoided heat", -2.000573183573656e-15], ["avoided electricity", -2.3674125085662277e-14]]'
$ws.Range("I3").Value = '[["CDU updated", 4.604976136023767e-09], ["Steel cylinder", 6.725000501277105e-11], ["nitrous oxide production", 2.7210021087733674e-08], ["Use phase CDU", 1.110894055900818e-08], ["EoL CDU", 2.9126087035565844e-10], ["EoL 50L cylinder - FU", 5.055227005369317e-09], ["avoided heat", -1.1298324241335168e-14], ["avoided electricity", -1.2012507594391547e-13]]'
$ws.Range("I4").Value = '[["MDU updated", 1.1984131040252541e-08], ["Steel cylinder", 6.725000501277105e-11], ["nitrous oxide production", 2.7210021087733674e-08], ["Use phase MDU", 5.774491740276279e-09], ["EoL MDU", 4.292758324859426e-09], ["EoL 50L cylinder - FU", 5.055227005369317e-09], ["avoided heat", -3.32544819537311e-14], ["avoided electricity", -5.78094907619593e-13]]'
$ws.Range("J3").Value = '[["CDU updated", 0.02297194643214596], ["Steel cylinder", 3.451402200181894e-05], ["nitrous oxide production", 0.42490983327819815], ["Use phase CDU", 0.0003478036811214869], ["EoL CDU", -0.0001129496014677351], ["EoL 50L cylinder - FU", 0.052218194516134374], ["avoided heat", -6.872013792021059e-08], ["avoided electricity", -3.1894403163726294e-08]]'
$ws.Range("J4").Value = '[["MDU updated", 0.05186765582766427], ["Steel cylinder", 3.451402200181894e-05], ["nitrous oxide production", 0.42490983327819815], ["Use phase MDU", 0.00013104305672162063], ["EoL MDU", -0.0005253249882754571], ["EoL 50L cylinder - FU", 0.052218194516134374], ["avoided heat", -2.0226473745363095e-07], ["avoided electricity", -1.5348995124963722e-07]]'
$ws.Range("K3").Value = '[["CDU updated", 1.4485294249440825], ["Steel cylinder", 0.05341058184255911], ["nitrous oxide production", 28.792690550160245], ["Use phase CDU", 30.062931903505202], ["EoL CDU", 0.028179963123186758], ["EoL 50L cylinder - FU", 4.5534051368332715], ["avoided heat", -3.023150091949239e-05], ["avoided electricity", -0.00010258268643979891]]'
$ws.Range("K4").Value = '[["MDU updated", 3.6900959202500014], ["Steel cylinder", 0.05341058184255911], ["nitrous oxide production", 28.792690550160245], ["Use phase MDU", 11.326902802431254], ["EoL MDU", 0.2979519424652686], ["EoL 50L cylinder - FU", 4.5534051368332715], ["avoided heat", -8.898070902776089e-05], ["avoided electricity", -0.0004936731832450832]]'
$ws.Range("L3").Value = '[["CDU updated", 2.3284601371953983e-05], ["Steel cylinder", -6.260798244751308e-08], ["nitrous oxide production", 3.528425210128439e-05], ["Use phase CDU", 9.73441197116697e-06], ["EoL CDU", -7.40860365471259e-07], ["EoL 50L cylinder - FU", 5.124795543471044e-06], ["avoided heat", -1.2015470710380128e-11], ["avoided electricity", -3.857211487781677e-11]]'
$ws.Range("L4").Value = '[["MDU updated", 7.762141349086539e-05], ["Steel cylinder", -6.260798244751308e-08], ["nitrous oxide production", 3.528425210128439e-05], ["Use phase MDU", 3.667664171704964e-06], ["EoL MDU", 8.076422627217764e-07], ["EoL 50L cylinder - FU", 5.124795543471044e-06], ["avoided heat", -3.536526704035473e-11], ["avoided electricity", -1.856260486145807e-10]]'
$ws.Range("M3").Value = '[["CDU updated", 6.4301405431153076e-09], ["Steel cylinder", 1.232676817653973e-10], ["nitrous oxide production", 8.530261057304866e-08], ["Use phase CDU", 1.9450074606640816e-09], ["EoL CDU", -1.3446117198089737e-11], ["EoL 50L cylinder - FU", 1.0432841391740725e-07], ["avoided heat", -9.197420007104196e-14], ["avoided electricity", -2.2874017791016e-13]]'
$ws.Range("M4").Value = '[["MDU updated", 1.8962760185204527e-08], ["Steel cylinder", 1.232676817653973e-10], ["nitrous oxide production", 8.530261057304866e-08], ["Use phase MDU", 7.328264098737524e-10], ["EoL MDU", -2.7492773949681775e-11], ["EoL 50L cylinder - FU", 1.0432841391740725e-07], ["avoided heat", -2.707086742390444e-13], ["avoided electricity", -1.1007987381487465e-12]]'
$ws.Range("N3").Value = '[["CDU updated", 2.9210859505135567e-08], ["Steel cylinder", 5.476745429210785e-10], ["nitrous oxide production", 5.905655038504484e-08], ["Use phase CDU", 1.3992546510312443e-08], ["EoL CDU", -7.652082336354675e-10], ["EoL 50L cylinder - FU", 1.6048194600218172e-08], ["avoided heat", -6.86255402237205e-14], ["avoided electricity", -1.323471443040989e-13]]'
$ws.Range("N4").Value = '[["MDU updated", 6.342169280464495e-08], ["Steel cylinder", 5.476745429210785e-10], ["nitrous oxide production", 5.905655038504484e-08], ["Use phase MDU", 5.272014545798531e-09], ["EoL MDU", -4.0508970365985256e-09], ["EoL 50L cylinder - FU", 1.6048194600218172e-08], ["avoided heat", -2.0198630701312296e-13], ["avoided electricity", -6.369128973197927e-13]]'
$ws.Range("O3").Value = '[["CDU updated", 0.0026802332411080343], ["Steel cylinder", 2.334134230290092e-05], ["nitrous oxide production", 0.007059262000778557], ["Use phase CDU", 0.000860498465395691], ["EoL CDU", -5.4689801265332504e-05], ["EoL 50L cylinder - FU", 0.0012673930563316849], ["avoided heat", -3.184894111340053e-09], ["avoided electricity", -1.9964517102311484e-08]]'
$ws.Range("O4").Value = '[["MDU updated", 0.005810953630333549], ["Steel cylinder", 2.334134230290092e-05], ["nitrous oxide production", 0.007059262000778557], ["Use phase MDU", 0.0003242126386973229], ["EoL MDU", -0.0002788564528856964], ["EoL 50L cylinder - FU", 0.0012673930563316849], ["avoided heat", -9.374133852787247e-09], ["avoided electricity", -9.607807178677165e-08]]'
$ws.Range("P3").Value = '[["CDU updated", 0.046205705935014954], ["Steel cylinder", 0.0015157174317935512], ["nitrous oxide production", 4.430496838846222], ["Use phase CDU", 0.026577952308677204], ["EoL CDU", -0.000902135155534033], ["EoL 50L cylinder - FU", 0.46994446269736545], ["avoided heat", -1.205177088932026e-07], ["avoided electricity", -9.317975012427849e-07]]'
$ws.Range("P4").Value = '[["MDU updated", 0.12369240376316787], ["Steel cylinder", 0.0015157174317935512], ["nitrous oxide production", 4.430496838846222], ["Use phase MDU", 0.010013856381725932], ["EoL MDU", -0.004574211602495358], ["EoL 50L cylinder - FU", 0.46994446269736545], ["avoided heat", -3.5472109756262925e-07], ["avoided electricity", -4.484221018537859e-06]]'
